# Applies the updated cryptocurrency price/volume data to the active worksheet.
# Most cells are plain text (coin names, URLs, percentage strings) and can be
# written directly. Some "Price" column values look like plain numbers (e.g.
# "1.00", "596.25") - writing those directly would make Excel silently convert
# them to numeric values, trimming formatting such as trailing zeros. To keep
# those cells as text (matching the source data), we temporarily force a text
# number format before assigning the value, then restore the original ("Normal")
# cell style so no visible formatting changes are introduced.

function Set-TextValue($ws, $cellRef, $text) {
    $ws.Range($cellRef).Value = $text
}

function Set-ForcedTextValue($ws, $cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates that are unambiguously text (names, links, percentages, or
# price strings containing more than one "." so Excel cannot parse them as a
# plain number).
$directUpdates = @(
    @("D2", "67.508.06"),
    @("E2", "  -4.12%  "),
    @("D3", "3.281.10"),
    @("E3", "  -6.07%  "),
    @("E4", "  -0.01%  "),
    @("E5", "  -3.62%  "),
    @("E6", "  -10.44%  "),
    @("E7", "  -0.22%  "),
    @("D8", "3.274.25"),
    @("E8", "  -6.24%  "),
    @("E9", "  -9.50%  "),
    @("E10", "  -11.92%  "),
    @("E11", "  -4.02%  "),
    @("E12", "  -11.05%  "),
    @("E13", "  -14.37%  "),
    @("E14", "  -9.14%  "),
    @("D15", "3.805.69"),
    @("E15", "  -6.26%  "),
    @("D16", "67.509.86"),
    @("E16", "  -4.23%  "),
    @("D17", "3.280.01"),
    @("E17", "  -6.31%  "),
    @("E18", "  -10.02%  "),
    @("E19", "  -5.95%  "),
    @("E20", "  -13.04%  "),
    @("E21", "  -12.77%  "),
    @("E22", "  -11.93%  "),
    @("E23", "  -11.44%  "),
    @("E24", "  -10.95%  "),
    @("E25", "  -11.42%  "),
    @("E26", "  -0.02%  "),
    @("E27", "  -10.67%  "),
    @("E28", "  -13.28%  "),
    @("B29", "RenderToken"),
    @("C29", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"),
    @("E29", "  -8.38%  "),
    @("B30", "EthereumClassic"),
    @("C30", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"),
    @("E30", "  -11.70%  "),
    @("E31", "  -6.34%  "),
    @("E32", "  -8.00%  "),
    @("E33", "  -17.15%  "),
    @("E34", "  -12.74%  "),
    @("E35", "  -10.05%  "),
    @("E36", "  -0.07%  "),
    @("E37", "  -6.96%  "),
    @("E38", "  -5.44%  "),
    @("E39", "  -11.84%  "),
    @("E40", "  -15.81%  "),
    @("E41", "  -10.45%  "),
    @("E42", "  -14.12%  "),
    @("D43", "2.954.89"),
    @("E43", "  -10.59%  "),
    @("E44", "  -11.15%  "),
    @("D45", "0.0₃0597"),
    @("E45", "  -15.81%  "),
    @("E46", "  -9.35%  "),
    @("E47", "  -14.10%  "),
    @("B48", "ThetaToken"),
    @("C48", "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"),
    @("E48", "  -15.29%  "),
    @("B49", "USDe"),
    @("C49", "https://coinranking.com/coin/exbfr2U-0+usde-usde"),
    @("E49", "  -0.03%  "),
    @("B50", "Stellar"),
    @("C50", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"),
    @("E50", "  -10.76%  "),
    @("B51", "Monero"),
    @("C51", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"),
    @("E51", "  -6.75%  ")
)

# Price-column updates whose new text looks like a plain decimal number and
# therefore need the text format forced so Excel keeps them as strings.
$forcedTextUpdates = @(
    @("D5", "596.25"),
    @("D6", "151.88"),
    @("D7", "0.999"),
    @("D9", "0.546"),
    @("D11", "6.88"),
    @("D13", "38.75"),
    @("D18", "538.00"),
    @("D20", "7.23"),
    @("D21", "15.20"),
    @("D22", "0.765"),
    @("D23", "7.90"),
    @("D24", "86.13"),
    @("D25", "13.64"),
    @("D26", "1.00"),
    @("D27", "3.28"),
    @("D28", "2.18"),
    @("D29", "8.09"),
    @("D30", "29.52"),
    @("D33", "6.66"),
    @("D34", "5.82"),
    @("D35", "534.85"),
    @("D37", "0.0458"),
    @("D38", "53.48"),
    @("D39", "0.0863"),
    @("D40", "9.08"),
    @("D41", "0.128"),
    @("D42", "2.82"),
    @("D44", "0.270"),
    @("D46", "2.21"),
    @("D47", "26.96"),
    @("D48", "2.36"),
    @("D49", "1.00"),
    @("D50", "0.114"),
    @("D51", "124.72")
)

foreach ($u in $directUpdates) {
    Set-TextValue $ws $u[0] $u[1]
}

foreach ($u in $forcedTextUpdates) {
    Set-ForcedTextValue $ws $u[0] $u[1]
}

